$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A width (new column inserted before former first data column) ---
$ws.Columns.Item(1).ColumnWidth = 12.5

# --- Header rows ---
# Row 3: add "Subject Code" label in A3 (mirrors D3, which already holds it)
$ws.Range("A3").Value = "Subject Code"

# --- Data rows 5-14: Subject No 1-10, same pattern as column B ---
$ws.Range("A5").Value = 1
$ws.Range("A6").Formula = "=A5+1"
$ws.Range("A7").Formula = "=A6+1"
$ws.Range("A8").Formula = "=A7+1"
$ws.Range("A9").Formula = "=A8+1"
$ws.Range("A10").Formula = "=A9+1"
$ws.Range("A11").Formula = "=A10+1"
$ws.Range("A12").Formula = "=A11+1"
$ws.Range("A13").Formula = "=A12+1"
$ws.Range("A14").Formula = "=A13+1"

# --- Data rows 15-28: Subject No restarts at 201, incrementing by 1 ---
$ws.Range("A15").Value = 201
$ws.Range("A16").Formula = "=A15+1"
$ws.Range("A17").Formula = "=A16+1"
$ws.Range("A18").Formula = "=A17+1"
$ws.Range("A19").Formula = "=A18+1"
$ws.Range("A20").Formula = "=A19+1"
$ws.Range("A21").Formula = "=A20+1"
$ws.Range("A22").Formula = "=A21+1"
$ws.Range("A23").Formula = "=A22+1"
$ws.Range("A24").Formula = "=A23+1"
$ws.Range("A25").Formula = "=A24+1"
$ws.Range("A26").Formula = "=A25+1"
$ws.Range("A27").Formula = "=A26+1"
$ws.Range("A28").Formula = "=A27+1"

# --- Update selection to C7 (also scrolls view back to A1) ---
$ws.Range("C7").Select()
